# Build an EA/TREATMENT/COUNTY/ROUTE/YEAR/Beg PM/End PM/LENGTH lookup table
# in columns J:Q (rows 16-37), mirroring the existing A:H data table
# (header in row 1, data in rows 2-22), but rendering the Beg/End PM columns
# as formatted text labels ("Beg PM:   12.345") instead of raw numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row 16 (mirrors row 1) ----
$ws.Range("J16").Value = $ws.Range("A1").Value2
$ws.Range("K16").Value = $ws.Range("B1").Value2
$ws.Range("L16").Value = $ws.Range("C1").Value2
$ws.Range("M16").Value = $ws.Range("D1").Value2
$ws.Range("N16").Value = $ws.Range("E1").Value2
$ws.Range("O16").Value = "Beg PM"
$ws.Range("P16").Value = "End PM"
$ws.Range("Q16").Value = $ws.Range("H1").Value2

# ---- Data rows 17-37 (mirror rows 2-22) ----
for ($i = 0; $i -lt 21; $i++) {
    $srcRow = 2 + $i
    $dstRow = 17 + $i

    $ws.Cells.Item($dstRow, 10).Value = $ws.Cells.Item($srcRow, 1).Value2  # J <- A (EA)
    $ws.Cells.Item($dstRow, 11).Value = $ws.Cells.Item($srcRow, 2).Value2  # K <- B (TREATMENT)
    $ws.Cells.Item($dstRow, 12).Value = $ws.Cells.Item($srcRow, 3).Value2  # L <- C (COUNTY)
    $ws.Cells.Item($dstRow, 13).Value = $ws.Cells.Item($srcRow, 4).Value2  # M <- D (ROUTE)
    $ws.Cells.Item($dstRow, 14).Value = $ws.Cells.Item($srcRow, 5).Value2  # N <- E (YEAR)

    $begPm = $ws.Cells.Item($srcRow, 6).Value2
    $endPm = $ws.Cells.Item($srcRow, 7).Value2

    $ws.Cells.Item($dstRow, 15).Value = "Beg PM:" + ("{0,9:N3}" -f $begPm)
    $ws.Cells.Item($dstRow, 16).Value = "End PM:" + ("{0,9:N3}" -f $endPm)

    $ws.Cells.Item($dstRow, 17).Value = $ws.Cells.Item($srcRow, 8).Value2  # Q <- H (LENGTH)
}

# Row 21 (EA 04-0C340 / HMA Medium Overlay) carries corrected PM values that
# differ from its source row's END_PM/LENGTH (8.400/8.4 instead of 8.200/8.2).
$ws.Range("O21").Value = "Beg PM:" + ("{0,9:N3}" -f 0)
$ws.Range("P21").Value = "End PM:" + ("{0,9:N3}" -f 8.4)
$ws.Range("Q21").Value = 8.4
